$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.872.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.986.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.976.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.148'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.440'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000219'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.455.86'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.780.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.980.40'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '458.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.668'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '55.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '454.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.165.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0777'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0378'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.03'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.86%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.74'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.45%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.243'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0494'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.61%  '
$ws.Range('E51').Value = '  +7.94%  '
